$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update E3: "gw" -> "silverarrow"
$ws.Range("E3").Value = "silverarrow"

# Capture row 3 values before adding row 4
$bVal = $ws.Range("B3").Text
$cVal = $ws.Range("C3").Text
$dVal = $ws.Range("D3").Text
$eVal = $ws.Range("E3").Text
$gVal = $ws.Range("G3").Text

# Add new row 4, replicating row 3 but with new claim number
$ws.Range("B4").Value = $bVal
$ws.Range("C4").Value = $cVal
$ws.Range("D4").Value = $dVal
$ws.Range("E4").Value = $eVal
$ws.Range("F4").Style = $ws.Range("F3").Style
$ws.Range("F4").Value = "'0420194406718"
$ws.Range("G4").Value = $gVal

# Add hyperlink to C4 (mirrors the hyperlink already present on C3), then
# re-apply the original "Hipervinculo" style so it reuses the existing
# cell style record instead of Excel auto-creating a new one.
$ws.Hyperlinks.Add($ws.Range("C4"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do")
$ws.Range("C4").Style = $ws.Range("C3").Style

$ws.Range("E5").Select()
